## THM v1.2, updated headers on all docs
##
## 1) Add five new rows to the checklist table for the "Battery Case" section.
## 2) Move the _GoBack bookmark from the trailing empty paragraph onto the
##    end of the last new row's text.
## 3) Update the header title, date, and revision text boxes (both the
##    drawing and VML-fallback copies).

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 0. Remove the old _GoBack bookmark (it gets re-created at the end of
#    the new "Assembly of battery case external insulation" row).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 1. Add the "Battery Case" rows to the first table.
# ---------------------------------------------------------------------
$t = $d.Tables(1)

function Set-EmptyCell($cell) {
    $cell.Range.InsertXML("<w:p $wns/>")
}

# Row: "Battery Case" (bold section header, manual hanging indent)
$row = $t.Rows.Add()
$xml = "<w:p $wns>" +
       "<w:pPr><w:ind w:left=""420"" w:hanging=""360""/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>" +
       "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Battery Case</w:t></w:r>" +
       "</w:p>"
$t.Cell($row.Index, 1).Range.InsertXML($xml)
Set-EmptyCell $t.Cell($row.Index, 2)
Set-EmptyCell $t.Cell($row.Index, 3)

# Row: "Internal layout of heaters and insulation into battery case / batteries"
$row = $t.Rows.Add()
$xml = "<w:p $wns>" +
       "<w:pPr><w:pStyle w:val=""ListParagraph""/></w:pPr>" +
       "<w:r><w:t>Internal layout of heaters and insulation into battery case / batteries</w:t></w:r>" +
       "</w:p>"
$t.Cell($row.Index, 1).Range.InsertXML($xml)
Set-EmptyCell $t.Cell($row.Index, 2)
Set-EmptyCell $t.Cell($row.Index, 3)

# Row: "Layout of external battery case insulation"
$row = $t.Rows.Add()
$xml = "<w:p $wns>" +
       "<w:pPr><w:pStyle w:val=""ListParagraph""/></w:pPr>" +
       "<w:r><w:t>Layout of external battery case insulation</w:t></w:r>" +
       "</w:p>"
$t.Cell($row.Index, 1).Range.InsertXML($xml)
Set-EmptyCell $t.Cell($row.Index, 2)
Set-EmptyCell $t.Cell($row.Index, 3)

# Row: "Assembly of battery case internals with insulation, heaters"
$row = $t.Rows.Add()
$xml = "<w:p $wns>" +
       "<w:pPr><w:pStyle w:val=""ListParagraph""/></w:pPr>" +
       "<w:r><w:t>Assembly of battery case internals with insulation, heaters</w:t></w:r>" +
       "</w:p>"
$t.Cell($row.Index, 1).Range.InsertXML($xml)
Set-EmptyCell $t.Cell($row.Index, 2)
Set-EmptyCell $t.Cell($row.Index, 3)

# Row: "Assembly of battery case external insulation" (carries the _GoBack bookmark)
$row = $t.Rows.Add()
$xml = "<w:p $wns>" +
       "<w:pPr><w:pStyle w:val=""ListParagraph""/></w:pPr>" +
       "<w:r><w:t>Assembly of battery case external insulation</w:t></w:r>" +
       "<w:bookmarkStart w:id=""0"" w:name=""_GoBack""/><w:bookmarkEnd w:id=""0""/>" +
       "</w:p>"
$t.Cell($row.Index, 1).Range.InsertXML($xml)
Set-EmptyCell $t.Cell($row.Index, 2)
Set-EmptyCell $t.Cell($row.Index, 3)

# ---------------------------------------------------------------------
# 2. Update header title text: "QUALIFICATION CHECKLIST" -> "THM QUAL CHECKLIST"
#    (split across two runs, same run formatting, in both header copies).
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "QUALIFICATION CHECKLIST"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "THM QUAL CHECKLIST"
$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2, $false, $false, $false, $false) | Out-Null

$sections = $d.Sections
for ($i = 1; $i -le $sections.Count; $i++) {
    $hdr = $sections.Item($i).Headers.Item(1)
    $hf = $hdr.Range.Find
    $hf.ClearFormatting()
    $hf.Text = "QUALIFICATION CHECKLIST"
    $hf.Replacement.ClearFormatting()
    $hf.Replacement.Text = "THM QUAL CHECKLIST"
    $hf.Execute($hf.Text, $false, $false, $false, $false, $false, $true, 1, $false, $hf.Replacement.Text, 2, $false, $false, $false, $false) | Out-Null
}

# ---------------------------------------------------------------------
# 3. Update the revision date: "2018-05-28" -> "2018-06-02"
# ---------------------------------------------------------------------
$sections = $d.Sections
for ($i = 1; $i -le $sections.Count; $i++) {
    $hdr = $sections.Item($i).Headers.Item(1)
    $hf = $hdr.Range.Find
    $hf.ClearFormatting()
    $hf.Text = "2018-05-28"
    $hf.Replacement.ClearFormatting()
    $hf.Replacement.Text = "2018-06-02"
    $hf.Execute($hf.Text, $false, $false, $false, $false, $false, $true, 1, $false, $hf.Replacement.Text, 2, $false, $false, $false, $false) | Out-Null
}

# ---------------------------------------------------------------------
# 4. Bump the revision number: "Revision: 1.1" -> "Revision: 1.2"
# ---------------------------------------------------------------------
$sections = $d.Sections
for ($i = 1; $i -le $sections.Count; $i++) {
    $hdr = $sections.Item($i).Headers.Item(1)
    $hf = $hdr.Range.Find
    $hf.ClearFormatting()
    $hf.Text = "Revision: 1.1"
    $hf.Replacement.ClearFormatting()
    $hf.Replacement.Text = "Revision: 1.2"
    $hf.Execute($hf.Text, $false, $false, $false, $false, $false, $true, 1, $false, $hf.Replacement.Text, 2, $false, $false, $false, $false) | Out-Null
}

Write-Output "done"
